# Partial-model sliding-window results (window 9): refresh predicted
# IPC values (column C) and dependent DELTA / DELTA^2 / TOTAL / MSE cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the original IPC-RO (actual) readings in column B get
# rewritten with floating point values that carry the same rounding as
# the freshly recomputed DELTA column (B = C - DELTA), so restate them too.
$bOverrides = @{
    5 = 30.03999999999999
    6 = 30.21000000000001
    20 = 31.65000000000001
    24 = 32.84999999999999
    25 = 32.90000000000001
    26 = 33.09999999999999
    27 = 33.40000000000001
    29 = 34.09999999999999
    30 = 34.40000000000001
    31 = 34.90000000000001
    37 = 37.90000000000001
    39 = 38.90000000000001
    40 = 39.40000000000001
    41 = 39.90000000000001
    42 = 40.09999999999999
    43 = 40.59999999999999
    44 = 40.90000000000001
    45 = 41.20000000000001
}

# New predicted IPC (column C) values for every data row (2-51).
$cValues = @{
    2 = 29.9100170135498
    3 = 29.88937187194824
    4 = 29.96999931335449
    5 = 30.16000175476074
    6 = 30.22909355163574
    7 = 30.33775901794434
    8 = 30.40864372253418
    9 = 30.61941337585449
    10 = 30.43139457702637
    11 = 30.40582847595215
    12 = 30.52097129821777
    13 = 30.62080192565918
    14 = 30.70599174499512
    15 = 31.09451866149902
    16 = 31.2536449432373
    17 = 31.33582878112793
    18 = 31.1854305267334
    19 = 31.38635063171387
    20 = 31.94420433044434
    21 = 32.52571487426758
    22 = 32.49223709106445
    23 = 32.56542587280273
    24 = 32.8115119934082
    25 = 32.99087142944336
    26 = 32.95927810668945
    27 = 33.53360366821289
    28 = 33.62392044067383
    29 = 33.83574295043945
    30 = 34.44326400756836
    31 = 35.08557510375977
    32 = 35.61227798461914
    33 = 35.89323043823242
    34 = 36.06398391723633
    35 = 36.57632064819336
    36 = 37.20608139038086
    37 = 37.85781478881836
    38 = 38.30781173706055
    39 = 39.00927352905273
    40 = 39.57085037231445
    41 = 39.59535598754883
    42 = 40.05672073364258
    43 = 40.49784469604492
    44 = 40.71113967895508
    45 = 41.16640853881836
    46 = 41.36802291870117
    47 = 42.05647659301758
    48 = 41.7586784362793
    49 = 43.38093185424805
    50 = 44.19462203979492
    51 = 43.60408782958984
}

$deltaSum = 0.0
$deltaSqSum = 0.0
$rowCount = 0

foreach ($row in 2..51) {
    if ($bOverrides.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value = $bOverrides[$row]
    }

    $b = $ws.Cells.Item($row, 2).Value()
    $c = $cValues[$row]
    $ws.Cells.Item($row, 3).Value = $c

    $delta = $c - $b
    $deltaSq = $delta * $delta
    $ws.Cells.Item($row, 4).Value = $delta
    $ws.Cells.Item($row, 5).Value = $deltaSq

    $deltaSum += $delta
    $deltaSqSum += $deltaSq
    $rowCount += 1
}

# Row 52 TOTAL: sum of DELTA and sum of DELTA^2 across all data rows.
$ws.Cells.Item(52, 3).Value = $deltaSum
$ws.Cells.Item(52, 5).Value = $deltaSqSum

# Row 53 MSE: mean of DELTA^2 across all data rows.
$ws.Cells.Item(53, 5).Value = $deltaSqSum / $rowCount

